$d = $word.ActiveDocument

# --- Change 1: etable(earnings_model) -> etable(earnings_model, vcov = 'hetero') ---
# Split the single "(earnings_model)" run into separately-styled runs that
# mirror the syntax-highlighted R call `(earnings_model, vcov = 'hetero')`.
$r1 = $d.Content
$r1.Find.Execute("(earnings_model)")
$r1.Text = "(earnings_model,"
$r1.Style = "NormalTok"

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("vcov =")
$r2.Style = "AttributeTok"

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(" ")
$r3.Style = "NormalTok"

$r4 = $d.Range($r3.End, $r3.End)
$r4.InsertAfter("'hetero'")
$r4.Style = "StringTok"

$r5 = $d.Range($r4.End, $r4.End)
$r5.InsertAfter(")")
$r5.Style = "NormalTok"

# --- Change 2: Constant std. error 0.0024 -> 0.0023 ---
$r = $d.Content
$r.Find.Execute("Constant                              11.44*** (0.0024)")
$r.Text = "Constant                              11.44*** (0.0023)"

# --- Change 3: age std. error 3.46e-5 -> 3.58e-5 ---
$r = $d.Content
$r.Find.Execute("age                                 0.0017*** (3.46e-5)")
$r.Text = "age                                 0.0017*** (3.58e-5)"

# --- Change 4: educ_factor Bachelor's degree std. error 0.0017 -> 0.0016 ---
# (Use direct Range.Text assignment rather than Find-Replace so the
# straight apostrophe in "Bachelor's" is not autocorrected to a curly quote.)
$r = $d.Content
$r.Find.Execute("educ_factorBachelor'sdegreeorhigher  0.3887*** (0.0017)")
$r.Text = "educ_factorBachelor'sdegreeorhigher  0.3887*** (0.0016)"

# --- Change 5: educ_factor High school graduate std. error 0.0017 -> 0.0016 ---
$r = $d.Content
$r.Find.Execute("educ_factorHighschoolgraduate       -0.2130*** (0.0017)")
$r.Text = "educ_factorHighschoolgraduate       -0.2130*** (0.0016)"

# --- Change 6: educ_factor No schooling/Primary std. error 0.0118 -> 0.0142 ---
$r = $d.Content
$r.Find.Execute("educ_factorNoschooling/Primary      -0.4217*** (0.0118)")
$r.Text = "educ_factorNoschooling/Primary      -0.4217*** (0.0142)"

# --- Change 7: educ_factor Some college std. error 0.0019 -> 0.0018 ---
$r = $d.Content
$r.Find.Execute("educ_factorSomecollege              -0.0504*** (0.0019)")
$r.Text = "educ_factorSomecollege              -0.0504*** (0.0018)"

# --- Change 8: educ_factor Some high school std. error 0.0076 -> 0.0080 ---
$r = $d.Content
$r.Find.Execute("educ_factorSomehighschool           -0.5174*** (0.0076)")
$r.Text = "educ_factorSomehighschool           -0.5174*** (0.0080)"

# --- Change 9: S.E. type IID -> Heteroskedast.-rob. ---
$r = $d.Content
$r.Find.Execute("S.E. type                                           IID")
$r.Text = "S.E. type                           Heteroskedast.-rob."

Write-Host "All changes applied"
